$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B2" = 0.2217898832684825
    "C2" = 0.5058365758754864
    "J2" = 0.01945525291828794
    "O2" = 0.003891050583657588
    "P2" = 0.1478599221789883
    "S2" = 0.1011673151750973
    "B3" = 0.02238805970149254
    "C3" = 0.02238805970149254
    "J3" = 0.03731343283582089
    "P3" = 0.7313432835820896
    "S3" = 0.1865671641791045
    "J4" = 0.0851063829787234
    "P4" = 0.5957446808510638
    "S4" = 0.3191489361702128
    "P5" = 0.5
    "S5" = 0.5
    "B6" = 0.05555555555555555
    "D6" = 0.007407407407407408
    "E6" = 0.003703703703703704
    "F6" = 0.1111111111111111
    "J6" = 0.2074074074074074
    "O6" = 0.01481481481481482
    "Q6" = 0.1333333333333333
    "R6" = 0.07037037037037037
    "S6" = 0.3962962962962963
    "B7" = 0.1162790697674419
    "D7" = 0.004651162790697674
    "F7" = 0.08837209302325581
    "J7" = 0.1302325581395349
    "O7" = 0.009302325581395349
    "Q7" = 0.1674418604651163
    "R7" = 0.08372093023255814
    "S7" = 0.4
    "B8" = 0.0826271186440678
    "D8" = 0.01059322033898305
    "F8" = 0.06567796610169492
    "J8" = 0.125
    "O8" = 0.006355932203389831
    "Q8" = 0.173728813559322
    "R8" = 0.06991525423728813
    "S8" = 0.4661016949152542
    "B9" = 0.08888888888888889
    "D9" = 0.005555555555555556
    "E9" = 0.005555555555555556
    "J9" = 0.1333333333333333
    "O9" = 0.01666666666666667
    "R9" = 0.07777777777777778
    "S9" = 0.4722222222222222
    "B10" = 0.08425531914893616
    "D10" = 0.03404255319148936
    "F10" = 0.0902127659574468
    "J10" = 0.1234042553191489
    "O10" = 0.01957446808510638
    "Q10" = 0.2059574468085106
    "R10" = 0.0748936170212766
    "S10" = 0.3676595744680851
    "G11" = 0.1395348837209302
    "J11" = 0.07641196013289037
    "K11" = 0.1926910299003322
    "L11" = 0.5681063122923588
    "S11" = 0.02325581395348837
    "G12" = 0.7630057803468208
    "J12" = 0.1907514450867052
    "K12" = 0.005780346820809248
    "L12" = 0.01734104046242774
    "S12" = 0.02312138728323699
    "G13" = 0.7704918032786885
    "J13" = 0.1967213114754098
    "S13" = 0.03278688524590164
    "F15" = 0.004347826086956522
    "H15" = 0.2
    "I15" = 0.0782608695652174
    "J15" = 0.3565217391304348
    "K15" = 0.07391304347826087
    "M15" = 0.01304347826086956
    "O15" = 0.06956521739130435
    "S15" = 0.2043478260869565
    "F16" = 0.03125
    "H16" = 0.23125
    "I16" = 0.0375
    "J16" = 0.35625
    "K16" = 0.1125
    "M16" = 0.01875
    "O16" = 0.05625
    "S16" = 0.15625
    "F17" = 0.01421800947867299
    "H17" = 0.2085308056872038
    "I17" = 0.06872037914691943
    "J17" = 0.3886255924170616
    "K17" = 0.0995260663507109
    "M17" = 0.01895734597156398
    "O17" = 0.07109004739336493
    "S17" = 0.1303317535545024
    "F18" = 0.02941176470588235
    "H18" = 0.1529411764705882
    "I18" = 0.09411764705882353
    "J18" = 0.3705882352941177
    "K18" = 0.1117647058823529
    "M18" = 0.04705882352941176
    "N18" = 0.005882352941176471
    "O18" = 0.09411764705882353
    "S18" = 0.09411764705882353
    "F19" = 0.01761102603369066
    "H19" = 0.217457886676876
    "I19" = 0.08652373660030628
    "J19" = 0.3330781010719755
    "K19" = 0.1094946401225115
    "M19" = 0.03215926493108729
    "N19" = 0.002297090352220521
    "O19" = 0.07044410413476264
    "S19" = 0.1309341500765697
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

Write-Output "Updated $($updates.Count) cells"
